# Update "By Vintage" sheet rankings with freshly re-scraped data (webscraper run 12/6).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("By Vintage")

# --- Row 21 (2006) ---
$ws.Range("B21").Value = 0.8181818181818182
$ws.Range("D21").Value = 0.1818181818181818

# --- Row 23 (2008) ---
$ws.Range("B23").Value = 0.7333333333333333
$ws.Range("D23").Value = 0.2333333333333333
$ws.Range("F23").Value = 0.03333333333333333
$ws.Range("G23").Value = 30

# --- Row 24 (2009) ---
$ws.Range("B24").Value = 0.6129032258064516
$ws.Range("D24").Value = 0.3548387096774194
$ws.Range("F24").Value = 0.03225806451612903
$ws.Range("G24").Value = 31

# --- Row 27 (2012) ---
$ws.Range("B27").Value = 0.6363636363636364
$ws.Range("D27").Value = 0.2727272727272727
$ws.Range("E27").Value = 0.0303030303030303
$ws.Range("F27").Value = 0.06060606060606061
$ws.Range("G27").Value = 33

# --- Row 29 (2014) ---
$ws.Range("B29").Value = 0.6024096385542169
$ws.Range("D29").Value = 0.2168674698795181
$ws.Range("E29").Value = 0.01204819277108434
$ws.Range("F29").Value = 0.1686746987951807
$ws.Range("G29").Value = 83

# --- Row 30 (2015) ---
$ws.Range("B30").Value = 0.5824175824175825
$ws.Range("C30").Value = 0.03296703296703297
$ws.Range("D30").Value = 0.2857142857142857
$ws.Range("F30").Value = 0.0989010989010989
$ws.Range("G30").Value = 91

# --- Row 31 (2016) ---
$ws.Range("B31").Value = 0.6131386861313869
$ws.Range("F31").Value = 0.218978102189781

# --- Row 32 (2017) ---
$ws.Range("B32").Value = 0.4951456310679612
$ws.Range("C32").Value = 0.01941747572815534
$ws.Range("D32").Value = 0.1067961165048544
$ws.Range("E32").Value = 0.009708737864077669
$ws.Range("F32").Value = 0.3689320388349515
$ws.Range("G32").Value = 206

# --- Row 33 (2018) ---
$ws.Range("B33").Value = 0.305
$ws.Range("C33").Value = 0.07000000000000001
$ws.Range("D33").Value = 0.05
$ws.Range("E33").Value = 0.07000000000000001
$ws.Range("F33").Value = 0.505
$ws.Range("G33").Value = 200

# --- Row 34 (2019) ---
$ws.Range("B34").Value = 0.2202797202797203
$ws.Range("C34").Value = 0.1048951048951049
$ws.Range("D34").Value = 0.01748251748251748
$ws.Range("E34").Value = 0.0944055944055944
$ws.Range("F34").Value = 0.5629370629370629
$ws.Range("G34").Value = 286

# --- Row 35 (2020) ---
$ws.Range("B35").Value = 0.1201044386422977
$ws.Range("C35").Value = 0.1462140992167102
$ws.Range("D35").Value = 0.01044386422976501
$ws.Range("E35").Value = 0.1044386422976501
$ws.Range("F35").Value = 0.618798955613577
$ws.Range("G35").Value = 383

# --- Row 36 (2021) ---
$ws.Range("B36").Value = 0.01937984496124031
$ws.Range("C36").Value = 0.08139534883720931
$ws.Range("D36").Value = 0.007751937984496124
$ws.Range("E36").Value = 0.07751937984496124
$ws.Range("F36").Value = 0.813953488372093
$ws.Range("G36").Value = 258

# --- New row 37 (2022) ---
$ws.Range("A37").Value = "2022"
$ws.Range("B37").Value = 0
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 1

# Match the formatting used by the other "Year Formed" label cells (bold, bordered, centered).
$ws.Range("A36").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$excel.CutCopyMode = $false
